$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 651
$ws.Range("C2").Value = 1776

$ws.Range("B3").Value = 645
$ws.Range("C3").Value = 383

$ws.Range("B4").Value = 401
$ws.Range("C4").Value = 145

$ws.Range("B5").Value = 11

$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 203
